$d = $word.ActiveDocument
$W = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Get-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) "Install PyGithub libary" -> split into runs with spell-check
#    proofErr markers around "PyGithub" and "libary".
# ---------------------------------------------------------------------
$idx1 = Get-ParaIndexByText $d "Install PyGithub libary"
$xml1 = @"
<w:p xmlns:w='$W'>
<w:r><w:t xml:space='preserve'>Install </w:t></w:r>
<w:proofErr w:type='spellStart'/>
<w:r><w:t>PyGithub</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
<w:r><w:t xml:space='preserve'> </w:t></w:r>
<w:proofErr w:type='spellStart'/>
<w:r><w:t>libary</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
</w:p>
"@
$null = $d.Paragraphs.Item($idx1).Range.InsertXML($xml1)

# ---------------------------------------------------------------------
# 2) "pip install pymongo" -> split into runs with spell-check proofErr
#    markers around "pymongo".
# ---------------------------------------------------------------------
$idx2 = Get-ParaIndexByText $d "pip install pymongo"
$xml2 = @"
<w:p xmlns:w='$W'>
<w:r><w:t xml:space='preserve'>pip install </w:t></w:r>
<w:proofErr w:type='spellStart'/>
<w:r><w:t>pymongo</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
</w:p>
"@
$null = $d.Paragraphs.Item($idx2).Range.InsertXML($xml2)

# ---------------------------------------------------------------------
# 3) Add the visualisation section. In the original document, the
#    paragraph that follows the docker paragraph is an empty paragraph,
#    followed by three more empty paragraphs. The middle three of those
#    four empty paragraphs become two new paragraphs:
#       "Visualisation stuff:"
#       "npm install recharts" (with proofErr around "npm")
#    This leaves one empty paragraph before and the remaining empty
#    paragraphs after, untouched.
# ---------------------------------------------------------------------
$idxDocker = Get-ParaIndexByText $d "Must have docker installed – for windows this means docker compose is not a separate installation."
if ($idxDocker -eq -1) {
    $idxDocker = Get-ParaIndexByText $d "Must have docker installed - for windows this means docker compose is not a separate installation."
}

$startPara = $d.Paragraphs.Item($idxDocker + 1)   # first empty paragraph after docker line
$endPara   = $d.Paragraphs.Item($idxDocker + 4)   # third empty paragraph following it

$r = $d.Range($startPara.Range.Start, $endPara.Range.End)
$xml3 = @"
<w:p xmlns:w='$W'><w:r><w:t>Visualisation stuff:</w:t></w:r></w:p>
<w:p xmlns:w='$W'><w:proofErr w:type='spellStart'/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> install recharts</w:t></w:r></w:p>
"@
$null = $r.InsertXML($xml3)
